# Edit script: applies the "24/12/2017 MAMATHA CHICK IN" changes.
$d = $word.ActiveDocument

# --- Change 1: merge the split "SUN Dec 17" / " 13:58:23 PST 2017" runs into one run ---
[void]$d.Content.Find.Execute("SUN Dec 17 13:58:23 PST 2017", $true, $false, $false, $false, $false, $true, 1, $false, "SUN Dec 17 13:58:23 PST 2017", 2)

# --- Change 2: append a new purchase-record block after the last "Amount Received mode ... - CASH" paragraph ---
$lastCash = $null
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $para = $d.Paragraphs($i)
    if ($para.Range.Text -like "*Amount Received mode*CASH*") {
        $lastCash = $para
        break
    }
}

$insertPoint = $d.Range($lastCash.Range.End, $lastCash.Range.End)
$newBlockXml = "<w:p><w:pPr><w:pStyle w:val=`"PlainText`"/><w:rPr><w:rFonts w:ascii=`"Courier New`" w:hAnsi=`"Courier New`" w:cs=`"Courier New`"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val=`"PlainText`"/><w:rPr><w:rFonts w:ascii=`"Courier New`" w:hAnsi=`"Courier New`" w:cs=`"Courier New`"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii=`"Courier New`" w:hAnsi=`"Courier New`" w:cs=`"Courier New`"/></w:rPr><w:t>SAT Dec 23</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=`"Courier New`" w:hAnsi=`"Courier New`" w:cs=`"Courier New`"/></w:rPr><w:t xml:space=`"preserve`"> 11:24:50 PST 2017</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=`"PlainText`"/><w:rPr><w:rFonts w:ascii=`"Courier New`" w:hAnsi=`"Courier New`" w:cs=`"Courier New`"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii=`"Courier New`" w:hAnsi=`"Courier New`" w:cs=`"Courier New`"/></w:rPr><w:t>Person Name</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=`"Courier New`" w:hAnsi=`"Courier New`" w:cs=`"Courier New`"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:rFonts w:ascii=`"Courier New`" w:hAnsi=`"Courier New`" w:cs=`"Courier New`"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:rFonts w:ascii=`"Courier New`" w:hAnsi=`"Courier New`" w:cs=`"Courier New`"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:rFonts w:ascii=`"Courier New`" w:hAnsi=`"Courier New`" w:cs=`"Courier New`"/></w:rPr><w:tab/><w:t>- GN</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=`"PlainText`"/><w:rPr><w:rFonts w:ascii=`"Courier New`" w:hAnsi=`"Courier New`" w:cs=`"Courier New`"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii=`"Courier New`" w:hAnsi=`"Courier New`" w:cs=`"Courier New`"/></w:rPr><w:t>Bill number</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=`"Courier New`" w:hAnsi=`"Courier New`" w:cs=`"Courier New`"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:rFonts w:ascii=`"Courier New`" w:hAnsi=`"Courier New`" w:cs=`"Courier New`"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:rFonts w:ascii=`"Courier New`" w:hAnsi=`"Courier New`" w:cs=`"Courier New`"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:rFonts w:ascii=`"Courier New`" w:hAnsi=`"Courier New`" w:cs=`"Courier New`"/></w:rPr><w:tab/><w:t>- 2125</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=`"PlainText`"/><w:rPr><w:rFonts w:ascii=`"Courier New`" w:hAnsi=`"Courier New`" w:cs=`"Courier New`"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii=`"Courier New`" w:hAnsi=`"Courier New`" w:cs=`"Courier New`"/></w:rPr><w:t>---------------------------------------------------------------</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=`"PlainText`"/><w:rPr><w:rFonts w:ascii=`"Courier New`" w:hAnsi=`"Courier New`" w:cs=`"Courier New`"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii=`"Courier New`" w:hAnsi=`"Courier New`" w:cs=`"Courier New`"/></w:rPr><w:t>Item Name</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=`"Courier New`" w:hAnsi=`"Courier New`" w:cs=`"Courier New`"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:rFonts w:ascii=`"Courier New`" w:hAnsi=`"Courier New`" w:cs=`"Courier New`"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:rFonts w:ascii=`"Courier New`" w:hAnsi=`"Courier New`" w:cs=`"Courier New`"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:rFonts w:ascii=`"Courier New`" w:hAnsi=`"Courier New`" w:cs=`"Courier New`"/></w:rPr><w:tab/><w:t>- CARROT2</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=`"PlainText`"/><w:rPr><w:rFonts w:ascii=`"Courier New`" w:hAnsi=`"Courier New`" w:cs=`"Courier New`"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii=`"Courier New`" w:hAnsi=`"Courier New`" w:cs=`"Courier New`"/></w:rPr><w:t>Number of Pockets</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=`"Courier New`" w:hAnsi=`"Courier New`" w:cs=`"Courier New`"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:rFonts w:ascii=`"Courier New`" w:hAnsi=`"Courier New`" w:cs=`"Courier New`"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:rFonts w:ascii=`"Courier New`" w:hAnsi=`"Courier New`" w:cs=`"Courier New`"/></w:rPr><w:tab/><w:t>- 1</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=`"PlainText`"/><w:rPr><w:rFonts w:ascii=`"Courier New`" w:hAnsi=`"Courier New`" w:cs=`"Courier New`"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii=`"Courier New`" w:hAnsi=`"Courier New`" w:cs=`"Courier New`"/></w:rPr><w:t>Number of KGs</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=`"Courier New`" w:hAnsi=`"Courier New`" w:cs=`"Courier New`"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:rFonts w:ascii=`"Courier New`" w:hAnsi=`"Courier New`" w:cs=`"Courier New`"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:rFonts w:ascii=`"Courier New`" w:hAnsi=`"Courier New`" w:cs=`"Courier New`"/></w:rPr><w:tab/><w:t>- 38</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=`"PlainText`"/><w:rPr><w:rFonts w:ascii=`"Courier New`" w:hAnsi=`"Courier New`" w:cs=`"Courier New`"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii=`"Courier New`" w:hAnsi=`"Courier New`" w:cs=`"Courier New`"/></w:rPr><w:t>Rate</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=`"Courier New`" w:hAnsi=`"Courier New`" w:cs=`"Courier New`"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:rFonts w:ascii=`"Courier New`" w:hAnsi=`"Courier New`" w:cs=`"Courier New`"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:rFonts w:ascii=`"Courier New`" w:hAnsi=`"Courier New`" w:cs=`"Courier New`"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:rFonts w:ascii=`"Courier New`" w:hAnsi=`"Courier New`" w:cs=`"Courier New`"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:rFonts w:ascii=`"Courier New`" w:hAnsi=`"Courier New`" w:cs=`"Courier New`"/></w:rPr><w:tab/><w:t>- 52</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=`"PlainText`"/><w:rPr><w:rFonts w:ascii=`"Courier New`" w:hAnsi=`"Courier New`" w:cs=`"Courier New`"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii=`"Courier New`" w:hAnsi=`"Courier New`" w:cs=`"Courier New`"/></w:rPr><w:t>Total Price</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=`"Courier New`" w:hAnsi=`"Courier New`" w:cs=`"Courier New`"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:rFonts w:ascii=`"Courier New`" w:hAnsi=`"Courier New`" w:cs=`"Courier New`"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:rFonts w:ascii=`"Courier New`" w:hAnsi=`"Courier New`" w:cs=`"Courier New`"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:rFonts w:ascii=`"Courier New`" w:hAnsi=`"Courier New`" w:cs=`"Courier New`"/></w:rPr><w:tab/><w:t>- 1976.0</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=`"PlainText`"/><w:rPr><w:rFonts w:ascii=`"Courier New`" w:hAnsi=`"Courier New`" w:cs=`"Courier New`"/><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii=`"Courier New`" w:hAnsi=`"Courier New`" w:cs=`"Courier New`"/><w:b/></w:rPr><w:t>Amount balance</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=`"Courier New`" w:hAnsi=`"Courier New`" w:cs=`"Courier New`"/><w:b/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:rFonts w:ascii=`"Courier New`" w:hAnsi=`"Courier New`" w:cs=`"Courier New`"/><w:b/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:rFonts w:ascii=`"Courier New`" w:hAnsi=`"Courier New`" w:cs=`"Courier New`"/><w:b/></w:rPr><w:tab/><w:t>- 7663.0</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=`"PlainText`"/><w:rPr><w:rFonts w:ascii=`"Courier New`" w:hAnsi=`"Courier New`" w:cs=`"Courier New`"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val=`"PlainText`"/><w:rPr><w:rFonts w:ascii=`"Courier New`" w:hAnsi=`"Courier New`" w:cs=`"Courier New`"/></w:rPr></w:pPr></w:p>"
[void]$insertPoint.InsertXML($newBlockXml)
